# "Minor updates and testing"
#
# 1. Fix spelling of a crew member's first name: "Stuart Charman" -> "Stewart
#    Charman" (appears three times, rows 46-48, column B).
#
# 2. The three "Joe Gunby" / Fireman rows (C19:C21) were stored as real date
#    serial values formatted with a custom "YYYY-MM-DD HH:MM:SS" number
#    format. Every other date in the sheet is stored as plain DD/MM/YYYY
#    text, so re-enter these three as text too, to match.
#
#    Typing "10/02/2018" straight into .Value would get reinterpreted by
#    Excel's smart date parser as a date serial (and, worse, misread as
#    2018-10-02 under a US locale) and forcing the cell's NumberFormat to
#    "@" to stop that permanently allocates a brand-new style record. To
#    avoid both problems, the literal text is produced with a throwaway
#    formula (a quoted string can't be reparsed as a date), copied, and
#    pasted back as a value - then any leftover number formatting on the
#    cell is cleared so it falls back to the sheet's default style, exactly
#    like the surrounding date cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the name typo (Grade=Trainee rows for this person) ---
$ws.Range("B46").Value = "Stewart Charman"
$ws.Range("B47").Value = "Stewart Charman"
$ws.Range("B48").Value = "Stewart Charman"

# --- 2. Re-enter the Fireman / Joe Gunby dates as plain text ---
function Set-DateAsText($cellAddr, $text) {
    $cell = $ws.Range($cellAddr)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)   # xlPasteValues
    $excel.CutCopyMode = $false
    $cell.ClearFormats()
}

Set-DateAsText "C19" "10/02/2018"
Set-DateAsText "C20" "11/02/2018"
Set-DateAsText "C21" "17/02/2018"

$wb.Save()
